# The edit replaces the entire letter body: the original two paragraphs
# (an "AI interests" paragraph and a "career goal" paragraph) are swapped
# out for a full cover-letter layout -- sender's address block, recipient's
# address block, salutation, three body paragraphs, sign-off and a typed
# signature -- while keeping the original "_GoBack" bookmark alive in its
# own paragraph.
#
# We build the replacement as a WordprocessingML fragment (a sequence of
# sibling <w:p> elements) and hand the whole thing to Range.InsertXML on
# $d.Content. InsertXML replaces the content of the range it is invoked on,
# so calling it on the full document Content swaps out every paragraph
# currently in the body for the ones constructed below, in one shot.

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function New-WRun($Text, $PreserveSpace) {
    # <w:r><w:t>...</w:t></w:r>, optionally with xml:space="preserve"
    $space = ""
    if ($PreserveSpace) { $space = " xml:space='preserve'" }
    return "<w:r><w:t$space>$Text</w:t></w:r>"
}

function New-WPara($InnerXml, $SpacingAfterZero) {
    # <w:p>, optionally carrying <w:pPr><w:spacing w:after="0"/></w:pPr>,
    # wrapping whatever run/bookmark XML is passed in.
    $pPr = ""
    if ($SpacingAfterZero) {
        $pPr = "<w:pPr><w:spacing w:after='0'/></w:pPr>"
    }
    return "<w:p xmlns:w='$wNs'>$pPr$InnerXml</w:p>"
}

function New-EmptyWPara() {
    return "<w:p xmlns:w='$wNs'/>"
}

$paragraphs = @()

# --- Sender's address block (each line has spacing-after = 0) ---
$paragraphs += New-WPara (New-WRun "Angad Singh" $false) $true
$paragraphs += New-WPara (New-WRun "59 Silkwood Cres " $true) $true
$paragraphs += New-WPara (New-WRun "Brampton, ON, Canada " $true) $true
$paragraphs += New-WPara (New-WRun "L6X 4K3" $false) $true
$paragraphs += New-WPara "" $true

# --- Recipient's address block (also spacing-after = 0) ---
$recipientLine1 = (New-WRun "Ottawa-Carleton" $false) + (New-WRun " Institute for Computer Science" $true)
$paragraphs += New-WPara $recipientLine1 $true
$paragraphs += New-WPara (New-WRun "800 King Edward, STE 1024" $false) $true
$recipientCityLine = (New-WRun "Ottawa, ON" $false) + (New-WRun ", Canada" $false)
$paragraphs += New-WPara $recipientCityLine $true
$paragraphs += New-WPara (New-WRun "K1N 6N5" $false) $true

# --- blank separator line before the salutation ---
$paragraphs += New-EmptyWPara

# --- Salutation ---
$paragraphs += New-WPara (New-WRun "To Whom It May Concern:" $false) $false

# --- Body paragraph 1: intro / background ---
$intro1 = "I am applying for the Computer Science Master’s Degree Program. I recently completed my Undergraduate Degree from the University of Toronto and came to have a passion for furthering my understanding and research in Computer Science.  "
$intro2 = "I found that I had an interest in artificial intelligence with the courses that I enjoyed taking. "
$intro3 = "These courses included Introduction into Machine Learning, Introduction into Artificial Intelligence and Introduction into Neural Networks. In my final year of University I also decided to do a research project in Optical Flow algorithms with my professor. In a group we created tools that could help other researchers create their own datasets to use to test their algorithms. "
$body1 = (New-WRun $intro1 $true) + (New-WRun $intro2 $true) + (New-WRun $intro3 $true)
$paragraphs += New-WPara $body1 $false

# --- Body paragraph 2: personal qualities / career goal ---
$qual1 = "I am a hardworking and a team player who is prepared to take tough challenges that I may face associated with the computer science program. "
$qual2 = "My goal is to find a career in a large company that is heavily investing the rese"
$qual3 = "arch in artificial intelligence and become a part of their software development team. "
$body2 = (New-WRun $qual1 $true) + (New-WRun $qual2 $false) + (New-WRun $qual3 $true)
$paragraphs += New-WPara $body2 $false

# --- Body paragraph 3: closing pitch ---
$close1 = "I believe that I have the capabilities and the perseverance to succeed the computer science master’s program. "
$close2 = "Mu strong desire to contribute to the development and growth of artificial intelligence research continues to drive me to succeed and to accomplish my goals. I have attached the rest of my application along with my referral letters. Thank you and I hope to speak you soon. "
$body3 = (New-WRun $close1 $true) + (New-WRun $close2 $true)
$paragraphs += New-WPara $body3 $false

# --- blank line, sign-off, blank signature line ---
$paragraphs += New-EmptyWPara
$paragraphs += New-WPara (New-WRun "Sincerely, " $true) $false
$paragraphs += New-EmptyWPara

# --- paragraph carrying the original _GoBack bookmark (kept, now alone) ---
$paragraphs += "<w:p xmlns:w='$wNs'><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"

# --- two trailing blank lines, then the typed signature ---
$paragraphs += New-EmptyWPara
$paragraphs += New-EmptyWPara
$paragraphs += New-WPara (New-WRun "Angad Singh " $true) $false

$newBodyXml = [string]::Join("", $paragraphs)

$d = $word.ActiveDocument
$d.Content.InsertXML($newBodyXml)

Write-Output "Paragraph count after edit: $($d.Paragraphs.Count)"
